$wb = $excel.ActiveWorkbook

# --- Update TxHash values on sheet "B2" ---
$wsB2 = $wb.Worksheets.Item("B2")
$wsB2.Range("A2").Value = "7B60F9194BA6D26830951F6DED4859EB050FB2891719BCEC108F1465D3F1F6C5"
$wsB2.Range("A3").Value = "F8F7D313F729B5562DB987A90A381CBDFE4FBFB33B45D370FC155D14316CE8DB"

# --- Update TxHash values on sheet "B1" ---
$wsB1 = $wb.Worksheets.Item("B1")
$wsB1.Range("A2").Value = "60CA2CBF6861A71DA1C24D5319A44BC8AEA544AD931E401234C3D004CE9F31D1"
$wsB1.Range("A3").Value = "164749B41D64B9A5D8081595AE3C075B4889F0CE678FFE596D33F26CB54A1169"

# --- Update selections / active sheet to match the saved workbook view state ---
# A20 keeps its old selection (E15) but loses tab focus - move focus away from it first.
$wsA20 = $wb.Worksheets.Item("A20")
[void]$wsA20.Activate()
[void]$wsA20.Range("E15").Select()

# B2 becomes selected at A4, but is not the final active tab.
[void]$wsB2.Activate()
[void]$wsB2.Range("A4").Select()

# B1 ends up as the active tab with H17 selected.
[void]$wsB1.Activate()
[void]$wsB1.Range("H17").Select()
